$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Extend row 2 with new values in C2:Y2, matching the formatting already
# used by A2/B2 (style index 1 -> Arial Unicode MS font).
$ws.Range("A2").Copy()
$ws.Range("C2:Y2").PasteSpecial(-4122)

$row2Values = @(3,6,7,6,5,5,8,6,5,5,7,6,7,8,7,7,7,5,7,7,5,6,7)
$col = 3
foreach ($val in $row2Values) {
    $ws.Cells.Item(2, $col).Value = $val
    $col = $col + 1
}

# Apply a yellow highlight fill to A25:B56 (creates new fill + cellXf entries)
$ws.Range("A25:B56").Interior.Color = 65535

# Update the active selection to match the committed workbook state
$null = $ws.Range("J13").Select()
